$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Cell B11 currently holds the text "R40" (shared string).
# Replace it with the text "1" (kept as text, not converted to the
# number 1, via the leading-apostrophe text-entry prefix) which becomes
# a new shared string entry.
$ws.Range("B11").Value = "'1"
